# Update the "b.md.md" row across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff has been generated ("Ready for handoff"),
# correcting the Latest Handoff File / Datetime for the b.md source file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is "b.md.md" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is "b.md.md" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-20 07:43:12"
foreach ($h in $zhcn.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 is "b.md.md" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-20 07:43:23"
foreach ($h in $dede.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
